# correção nos dados e inicio da analise PNAD 2009
#
# The sheet had a stray label-only row ("grandes regiões e unidades da
# federação") at row 6 that carried no data, which had pushed every
# region's values down by one row relative to its label (e.g. the row
# labeled "norte" actually held "rondônia"'s numbers, and so on all the
# way down, with "goiás"'s numbers trailing off in an extra row 37).
#
# Fix: remove that stray row entirely so labels and data line up again;
# Excel shifts every row below it up by one, which also drops the now
# unused shared string from the workbook when it is saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
